$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update response-option labels for the "inverse" (positively-worded, reverse-scored)
# SDQ items: 0-2 scale shown with the order flipped (2=Not True ... 0=Certainly True)
$newLabel = "2=Not True, 1=Somewhat True, 0=Certainly True"

$ws.Range("E3").Value = $newLabel
$ws.Range("E6").Value = $newLabel
$ws.Range("E11").Value = $newLabel
$ws.Range("E19").Value = $newLabel
$ws.Range("E22").Value = $newLabel

# Row 35 (SDQ_30): value range relabeled from 0-2 to 0-3
$ws.Range("D35").Value = "0-3"

# Column A resized (best-fit) to fit the longer question text
$ws.Columns.Item(1).ColumnWidth = 140.66666666666666

# Selection moved to E22 (last touched cell)
$ws.Range("E22").Select()
